$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = $null
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = $null
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = $null
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = $null
$ws.Range("O2").Value = $null
$ws.Range("P2").Value = $null
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = $null
$ws.Range("S2").Value = $null
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null
$ws.Range("G3").Value = $null
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = $null
$ws.Range("K3").Value = $null
$ws.Range("L3").Value = $null
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = $null
$ws.Range("P3").Value = $null
$ws.Range("Q3").Value = $null
$ws.Range("R3").Value = $null
$ws.Range("S3").Value = $null
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = $null
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = $null
$ws.Range("K4").Value = $null
$ws.Range("L4").Value = $null
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = $null
$ws.Range("O4").Value = $null
$ws.Range("P4").Value = $null
$ws.Range("Q4").Value = $null
$ws.Range("R4").Value = $null
$ws.Range("S4").Value = 2
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = $null
$ws.Range("H5").Value = $null
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = $null
$ws.Range("K5").Value = $null
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = $null
$ws.Range("O5").Value = $null
$ws.Range("P5").Value = $null
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = $null
$ws.Range("S5").Value = $null
$ws.Range("B6").Value = $null
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = $null
$ws.Range("F6").Value = $null
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = $null
$ws.Range("I6").Value = $null
$ws.Range("J6").Value = $null
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = $null
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = $null
$ws.Range("O6").Value = $null
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = $null
$ws.Range("R6").Value = $null
$ws.Range("S6").Value = $null
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null
$ws.Range("D7").Value = $null
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = $null
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = $null
$ws.Range("K7").Value = $null
$ws.Range("L7").Value = $null
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null
$ws.Range("O7").Value = $null
$ws.Range("P7").Value = $null
$ws.Range("Q7").Value = $null
$ws.Range("R7").Value = $null
$ws.Range("S7").Value = $null
$ws.Range("B8").Value = $null
$ws.Range("C8").Value = $null
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = $null
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = $null
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = $null
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = $null
$ws.Range("L8").Value = $null
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = $null
$ws.Range("O8").Value = $null
$ws.Range("P8").Value = $null
$ws.Range("Q8").Value = $null
$ws.Range("R8").Value = $null
$ws.Range("S8").Value = $null
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = $null
$ws.Range("E9").Value = $null
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $null
$ws.Range("H9").Value = $null
$ws.Range("I9").Value = $null
$ws.Range("J9").Value = $null
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = $null
$ws.Range("M9").Value = $null
$ws.Range("N9").Value = $null
$ws.Range("O9").Value = 3
$ws.Range("P9").Value = $null
$ws.Range("Q9").Value = $null
$ws.Range("R9").Value = $null
$ws.Range("S9").Value = $null
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = $null
$ws.Range("F10").Value = $null
$ws.Range("G10").Value = $null
$ws.Range("H10").Value = $null
$ws.Range("I10").Value = $null
$ws.Range("J10").Value = $null
$ws.Range("K10").Value = $null
$ws.Range("L10").Value = $null
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = $null
$ws.Range("O10").Value = $null
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = $null
$ws.Range("R10").Value = $null
$ws.Range("S10").Value = 1
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = $null
$ws.Range("F11").Value = $null
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = $null
$ws.Range("I11").Value = $null
$ws.Range("J11").Value = $null
$ws.Range("K11").Value = $null
$ws.Range("L11").Value = $null
$ws.Range("M11").Value = $null
$ws.Range("N11").Value = $null
$ws.Range("O11").Value = 3
$ws.Range("P11").Value = 2
$ws.Range("Q11").Value = $null
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = $null
$ws.Range("B12").Value = $null
$ws.Range("C12").Value = $null
$ws.Range("D12").Value = $null
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = $null
$ws.Range("G12").Value = $null
$ws.Range("H12").Value = $null
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = $null
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = $null
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = $null
$ws.Range("O12").Value = $null
$ws.Range("P12").Value = $null
$ws.Range("Q12").Value = $null
$ws.Range("R12").Value = $null
$ws.Range("S12").Value = $null
$ws.Range("B13").Value = $null
$ws.Range("C13").Value = $null
$ws.Range("D13").Value = $null
$ws.Range("E13").Value = $null
$ws.Range("F13").Value = $null
$ws.Range("G13").Value = $null
$ws.Range("H13").Value = $null
$ws.Range("I13").Value = $null
$ws.Range("J13").Value = $null
$ws.Range("K13").Value = $null
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = 2
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = $null
$ws.Range("Q13").Value = $null
$ws.Range("R13").Value = $null
$ws.Range("S13").Value = $null
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = $null
$ws.Range("D14").Value = $null
$ws.Range("E14").Value = $null
$ws.Range("F14").Value = $null
$ws.Range("G14").Value = $null
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = $null
$ws.Range("J14").Value = $null
$ws.Range("K14").Value = $null
$ws.Range("L14").Value = $null
$ws.Range("M14").Value = $null
$ws.Range("N14").Value = $null
$ws.Range("O14").Value = $null
$ws.Range("P14").Value = 1
$ws.Range("Q14").Value = $null
$ws.Range("R14").Value = $null
$ws.Range("S14").Value = $null
$ws.Range("B15").Value = $null
$ws.Range("C15").Value = $null
$ws.Range("D15").Value = $null
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = $null
$ws.Range("G15").Value = $null
$ws.Range("H15").Value = $null
$ws.Range("I15").Value = $null
$ws.Range("J15").Value = $null
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = $null
$ws.Range("M15").Value = $null
$ws.Range("N15").Value = $null
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = $null
$ws.Range("Q15").Value = $null
$ws.Range("R15").Value = $null
$ws.Range("S15").Value = $null
$ws.Range("B16").Value = $null
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = $null
$ws.Range("E16").Value = $null
$ws.Range("F16").Value = $null
$ws.Range("G16").Value = $null
$ws.Range("H16").Value = $null
$ws.Range("I16").Value = $null
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = $null
$ws.Range("L16").Value = $null
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = $null
$ws.Range("O16").Value = $null
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = $null
$ws.Range("R16").Value = $null
$ws.Range("S16").Value = $null
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = $null
$ws.Range("E17").Value = $null
$ws.Range("F17").Value = $null
$ws.Range("G17").Value = $null
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = $null
$ws.Range("J17").Value = $null
$ws.Range("K17").Value = $null
$ws.Range("L17").Value = $null
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = $null
$ws.Range("O17").Value = $null
$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = $null
$ws.Range("R17").Value = $null
$ws.Range("S17").Value = $null
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = $null
$ws.Range("D18").Value = $null
$ws.Range("E18").Value = $null
$ws.Range("F18").Value = $null
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = $null
$ws.Range("I18").Value = $null
$ws.Range("J18").Value = $null
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = $null
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = $null
$ws.Range("O18").Value = $null
$ws.Range("P18").Value = $null
$ws.Range("Q18").Value = $null
$ws.Range("R18").Value = $null
$ws.Range("S18").Value = $null
$ws.Range("B19").Value = $null
$ws.Range("C19").Value = $null
$ws.Range("D19").Value = $null
$ws.Range("E19").Value = $null
$ws.Range("F19").Value = $null
$ws.Range("G19").Value = $null
$ws.Range("H19").Value = $null
$ws.Range("I19").Value = $null
$ws.Range("J19").Value = $null
$ws.Range("K19").Value = $null
$ws.Range("L19").Value = 3
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = $null
$ws.Range("Q19").Value = $null
$ws.Range("R19").Value = $null
$ws.Range("S19").Value = $null
$ws.Range("B20").Value = $null
$ws.Range("C20").Value = $null
$ws.Range("D20").Value = $null
$ws.Range("E20").Value = $null
$ws.Range("F20").Value = $null
$ws.Range("G20").Value = $null
$ws.Range("H20").Value = $null
$ws.Range("I20").Value = $null
$ws.Range("J20").Value = $null
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = $null
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = $null
$ws.Range("O20").Value = 3
$ws.Range("P20").Value = $null
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = $null
$ws.Range("S20").Value = $null
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = $null
$ws.Range("D21").Value = $null
$ws.Range("E21").Value = $null
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = $null
$ws.Range("H21").Value = $null
$ws.Range("I21").Value = $null
$ws.Range("J21").Value = $null
$ws.Range("K21").Value = $null
$ws.Range("L21").Value = $null
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = $null
$ws.Range("O21").Value = $null
$ws.Range("P21").Value = $null
$ws.Range("Q21").Value = $null
$ws.Range("R21").Value = $null
$ws.Range("S21").Value = $null
$ws.Range("B22").Value = $null
$ws.Range("C22").Value = $null
$ws.Range("D22").Value = $null
$ws.Range("E22").Value = $null
$ws.Range("F22").Value = $null
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = $null
$ws.Range("I22").Value = $null
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = $null
$ws.Range("L22").Value = $null
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = $null
$ws.Range("O22").Value = $null
$ws.Range("P22").Value = $null
$ws.Range("Q22").Value = $null
$ws.Range("R22").Value = $null
$ws.Range("S22").Value = 2
$ws.Range("B23").Value = $null
$ws.Range("C23").Value = $null
$ws.Range("D23").Value = $null
$ws.Range("E23").Value = $null
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = $null
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = $null
$ws.Range("K23").Value = $null
$ws.Range("L23").Value = $null
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = $null
$ws.Range("O23").Value = $null
$ws.Range("P23").Value = $null
$ws.Range("Q23").Value = $null
$ws.Range("R23").Value = $null
$ws.Range("S23").Value = $null
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = $null
$ws.Range("D24").Value = $null
$ws.Range("E24").Value = $null
$ws.Range("F24").Value = $null
$ws.Range("G24").Value = $null
$ws.Range("H24").Value = $null
$ws.Range("I24").Value = $null
$ws.Range("J24").Value = $null
$ws.Range("K24").Value = $null
$ws.Range("L24").Value = 2
$ws.Range("M24").Value = $null
$ws.Range("N24").Value = $null
$ws.Range("O24").Value = $null
$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = $null
$ws.Range("R24").Value = $null
$ws.Range("S24").Value = $null
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = $null
$ws.Range("D25").Value = $null
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = $null
$ws.Range("G25").Value = $null
$ws.Range("H25").Value = $null
$ws.Range("I25").Value = $null
$ws.Range("J25").Value = $null
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = $null
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = $null
$ws.Range("O25").Value = $null
$ws.Range("P25").Value = $null
$ws.Range("Q25").Value = $null
$ws.Range("R25").Value = $null
$ws.Range("S25").Value = $null
$ws.Range("B26").Value = $null
$ws.Range("C26").Value = $null
$ws.Range("D26").Value = $null
$ws.Range("E26").Value = $null
$ws.Range("F26").Value = $null
$ws.Range("G26").Value = $null
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = $null
$ws.Range("J26").Value = $null
$ws.Range("K26").Value = $null
$ws.Range("L26").Value = $null
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = $null
$ws.Range("O26").Value = 2
$ws.Range("P26").Value = $null
$ws.Range("Q26").Value = $null
$ws.Range("R26").Value = $null
$ws.Range("S26").Value = 3
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = $null
$ws.Range("D27").Value = $null
$ws.Range("E27").Value = $null
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $null
$ws.Range("H27").Value = $null
$ws.Range("I27").Value = $null
$ws.Range("J27").Value = $null
$ws.Range("K27").Value = $null
$ws.Range("L27").Value = $null
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = 2
$ws.Range("O27").Value = $null
$ws.Range("P27").Value = $null
$ws.Range("Q27").Value = $null
$ws.Range("R27").Value = $null
$ws.Range("S27").Value = $null
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = $null
$ws.Range("D28").Value = $null
$ws.Range("E28").Value = $null
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = $null
$ws.Range("H28").Value = $null
$ws.Range("I28").Value = $null
$ws.Range("J28").Value = $null
$ws.Range("K28").Value = 1
$ws.Range("L28").Value = $null
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = $null
$ws.Range("O28").Value = $null
$ws.Range("P28").Value = $null
$ws.Range("Q28").Value = $null
$ws.Range("R28").Value = $null
$ws.Range("S28").Value = $null
$ws.Range("B29").Value = $null
$ws.Range("C29").Value = $null
$ws.Range("D29").Value = $null
$ws.Range("E29").Value = 3
$ws.Range("F29").Value = $null
$ws.Range("G29").Value = $null
$ws.Range("H29").Value = $null
$ws.Range("I29").Value = $null
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = $null
$ws.Range("L29").Value = $null
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = $null
$ws.Range("O29").Value = 1
$ws.Range("P29").Value = $null
$ws.Range("Q29").Value = $null
$ws.Range("R29").Value = $null
$ws.Range("S29").Value = $null
$ws.Range("B30").Value = $null
$ws.Range("C30").Value = $null
$ws.Range("D30").Value = $null
$ws.Range("E30").Value = $null
$ws.Range("F30").Value = $null
$ws.Range("G30").Value = $null
$ws.Range("H30").Value = $null
$ws.Range("I30").Value = $null
$ws.Range("J30").Value = $null
$ws.Range("K30").Value = $null
$ws.Range("L30").Value = $null
$ws.Range("M30").Value = $null
$ws.Range("N30").Value = $null
$ws.Range("O30").Value = $null
$ws.Range("P30").Value = 3
$ws.Range("Q30").Value = $null
$ws.Range("R30").Value = 2
$ws.Range("S30").Value = 1
$ws.Range("B31").Value = $null
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = $null
$ws.Range("E31").Value = $null
$ws.Range("F31").Value = $null
$ws.Range("G31").Value = $null
$ws.Range("H31").Value = $null
$ws.Range("I31").Value = $null
$ws.Range("J31").Value = $null
$ws.Range("K31").Value = $null
$ws.Range("L31").Value = $null
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = $null
$ws.Range("O31").Value = 2
$ws.Range("P31").Value = $null
$ws.Range("Q31").Value = $null
$ws.Range("R31").Value = $null
$ws.Range("S31").Value = 3
$ws.Range("B32").Value = $null
$ws.Range("C32").Value = $null
$ws.Range("D32").Value = $null
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = $null
$ws.Range("G32").Value = $null
$ws.Range("H32").Value = $null
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = $null
$ws.Range("K32").Value = $null
$ws.Range("L32").Value = $null
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = $null
$ws.Range("O32").Value = $null
$ws.Range("P32").Value = $null
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = $null
$ws.Range("S32").Value = $null
$ws.Range("B33").Value = $null
$ws.Range("C33").Value = $null
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = $null
$ws.Range("F33").Value = $null
$ws.Range("G33").Value = $null
$ws.Range("H33").Value = $null
$ws.Range("I33").Value = $null
$ws.Range("J33").Value = 3
$ws.Range("K33").Value = $null
$ws.Range("L33").Value = $null
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = $null
$ws.Range("P33").Value = $null
$ws.Range("Q33").Value = $null
$ws.Range("R33").Value = $null
$ws.Range("S33").Value = $null
$ws.Range("B34").Value = $null
$ws.Range("C34").Value = $null
$ws.Range("D34").Value = $null
$ws.Range("E34").Value = $null
$ws.Range("F34").Value = $null
$ws.Range("G34").Value = $null
$ws.Range("H34").Value = $null
$ws.Range("I34").Value = 3
$ws.Range("J34").Value = 2
$ws.Range("K34").Value = $null
$ws.Range("L34").Value = $null
$ws.Range("M34").Value = 1
$ws.Range("N34").Value = $null
$ws.Range("O34").Value = $null
$ws.Range("P34").Value = $null
$ws.Range("Q34").Value = $null
$ws.Range("R34").Value = $null
$ws.Range("S34").Value = $null
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = $null
$ws.Range("D35").Value = $null
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = $null
$ws.Range("G35").Value = $null
$ws.Range("H35").Value = $null
$ws.Range("I35").Value = $null
$ws.Range("J35").Value = $null
$ws.Range("K35").Value = $null
$ws.Range("L35").Value = $null
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = $null
$ws.Range("O35").Value = $null
$ws.Range("P35").Value = 2
$ws.Range("Q35").Value = $null
$ws.Range("R35").Value = $null
$ws.Range("S35").Value = $null
$ws.Range("B36").Value = $null
$ws.Range("C36").Value = $null
$ws.Range("D36").Value = $null
$ws.Range("E36").Value = $null
$ws.Range("F36").Value = 3
$ws.Range("G36").Value = $null
$ws.Range("H36").Value = $null
$ws.Range("I36").Value = $null
$ws.Range("J36").Value = 2
$ws.Range("K36").Value = $null
$ws.Range("L36").Value = $null
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = $null
$ws.Range("O36").Value = 1
$ws.Range("P36").Value = $null
$ws.Range("Q36").Value = $null
$ws.Range("R36").Value = $null
$ws.Range("S36").Value = $null
$ws.Range("B37").Value = 3
$ws.Range("C37").Value = $null
$ws.Range("D37").Value = $null
$ws.Range("E37").Value = $null
$ws.Range("F37").Value = $null
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = $null
$ws.Range("I37").Value = $null
$ws.Range("J37").Value = 1
$ws.Range("K37").Value = $null
$ws.Range("L37").Value = $null
$ws.Range("M37").Value = $null
$ws.Range("N37").Value = $null
$ws.Range("O37").Value = $null
$ws.Range("P37").Value = $null
$ws.Range("Q37").Value = $null
$ws.Range("R37").Value = $null
$ws.Range("S37").Value = $null

$ws.Range("P27").Select()
